# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.335.10"
$ws.Range("E2").Value = "  +0.55%  "

$ws.Range("D3").Value = "'1.880.15"
$ws.Range("E3").Value = "  +4.08%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'312.66"
$ws.Range("E5").Value = "  -0.19%  "

$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.00%  "

$ws.Range("D7").Value = "'0.5019"
$ws.Range("E7").Value = "  -2.15%  "

$ws.Range("D8").Value = "'0.3962"
$ws.Range("E8").Value = "  +0.31%  "

$ws.Range("D9").Value = "'0.09870"
$ws.Range("E9").Value = "  +26.70%  "

$ws.Range("D10").Value = "'1.128"
$ws.Range("E10").Value = "  +1.95%  "

$ws.Range("D11").Value = "'41.42"
$ws.Range("E11").Value = "  +1.23%  "

$ws.Range("D12").Value = "'6.484"
$ws.Range("E12").Value = "  +1.75%  "

$ws.Range("D13").Value = "'21.04"
$ws.Range("E13").Value = "  +3.20%  "

$ws.Range("D14").Value = "'1.867.35"
$ws.Range("E14").Value = "  +3.53%  "

$ws.Range("D15").Value = "'1.001"
$ws.Range("E15").Value = "  -0.03%  "

$ws.Range("D16").Value = "'7.401"
$ws.Range("E16").Value = "  +1.27%  "

$ws.Range("D17").Value = "'0.00001142"
$ws.Range("E17").Value = "  +6.06%  "

$ws.Range("D18").Value = "'93.74"
$ws.Range("E18").Value = "  +0.90%  "

$ws.Range("D19").Value = "'0.06682"
$ws.Range("E19").Value = "  +1.66%  "

$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "'17.41"
$ws.Range("E20").Value = "  +0.84%  "

$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "  -0.04%  "

$ws.Range("D22").Value = "'6.106"
$ws.Range("E22").Value = "  +1.57%  "

$ws.Range("D23").Value = "'28.393.08"
$ws.Range("E23").Value = "  +0.61%  "

$ws.Range("E24").Value = "  +2.13%  "

$ws.Range("D25").Value = "'2.260"
$ws.Range("E25").Value = "  +2.04%  "

$ws.Range("D26").Value = "'2.570"
$ws.Range("E26").Value = "  +5.31%  "

$ws.Range("B27").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C27").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D27").Value = "'2.090.39"
$ws.Range("E27").Value = "  +3.84%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'21.25"
$ws.Range("E28").Value = "  +3.68%  "

$ws.Range("D29").Value = "'158.04"
$ws.Range("E29").Value = "  -1.58%  "

$ws.Range("D30").Value = "'127.74"
$ws.Range("E30").Value = "  -0.17%  "

$ws.Range("E31").Value = "  -2.99%  "

$ws.Range("D32").Value = "'1.064"
$ws.Range("E32").Value = "  +0.64%  "

$ws.Range("D33").Value = "'5.649"
$ws.Range("E33").Value = "  +1.59%  "

$ws.Range("D34").Value = "'3.610"
$ws.Range("E34").Value = "  -1.35%  "

$ws.Range("D35").Value = "'0.06822"
$ws.Range("E35").Value = "  -4.47%  "

$ws.Range("D36").Value = "'9.471"
$ws.Range("E36").Value = "  +3.23%  "

$ws.Range("D37").Value = "'0.02393"
$ws.Range("E37").Value = "  +1.92%  "

$ws.Range("E38").Value = "  +1.07%  "

$ws.Range("D39").Value = "'5.035"
$ws.Range("E39").Value = "  -0.03%  "

$ws.Range("D40").Value = "'11.52"
$ws.Range("E40").Value = "  +0.03%  "

$ws.Range("D41").Value = "'0.6312"
$ws.Range("E41").Value = "  +2.55%  "

$ws.Range("D42").Value = "'1.177"
$ws.Range("E42").Value = "  +2.28%  "

$ws.Range("E43").Value = "  -0.03%  "

$ws.Range("D44").Value = "'13.48"
$ws.Range("E44").Value = "  +2.75%  "

$ws.Range("D45").Value = "'0.6019"
$ws.Range("E45").Value = "  +1.14%  "

$ws.Range("D46").Value = "'1.280"
$ws.Range("E46").Value = "  -1.74%  "

$ws.Range("D47").Value = "'3.683"
$ws.Range("E47").Value = "  -1.39%  "

$ws.Range("D48").Value = "'125.44"
$ws.Range("E48").Value = "  +0.60%  "

$ws.Range("D49").Value = "'1.995"
$ws.Range("E49").Value = "  +4.22%  "

$ws.Range("D50").Value = "'1.202"
$ws.Range("E50").Value = "  -0.77%  "

$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").Value = "'1.130"
$ws.Range("E51").Value = "  +5.70%  "
